# Bug fix: the English source string "Initial & Final Surveillance
# Diagnosis" used a literal "&" which broke the Lao translation lookup.
# The fix splits this into its own row with the corrected English text
# ("&" -> "and"), so a new row must be inserted and the old Lao
# translation (previously sitting next to the "&" string) is carried
# over to the new row, leaving the original row marked "TBT" (to be
# translated) since its translation moved away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 75; rows 75..181 shift down to 76..182.
$ws.Rows.Item(75).Insert()

# Row 74 keeps its English text ("Initial & Final Surveillance
# Diagnosis") but loses its Lao translation, which moves to row 75.
$ws.Cells.Item(74, 2).Value = "TBT"

# New row 75: corrected English text + the Lao translation that used
# to belong to row 74.
$ws.Cells.Item(75, 1).Value = "Initial and Final Surveillance Diagnosis"
$ws.Cells.Item(75, 2).Value = "ການບົ່ງມະຕິການເຝົ້າລະວັງໃນເບື້ອງຕົ້ນແລະຂັ້ນສຸດທ້າຍ"

# Same "&" -> "and" translation bug fix for the other affected string,
# which after the row insertion above now lives at row 145.
$ws.Cells.Item(145, 1).Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
